# Update cryptos list values (prices and 1h volume %) to match the
# latest scrape. A handful of rows also had their coin re-ranked,
# so Coin/Link (columns B/C) shifted along with Price/Volume (D/E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking price strings (e.g. "0.9995", "1.000", "7.540") must keep
# their exact text representation (leading/trailing zeros, no type change),
# so they are written with a leading apostrophe to force literal text -
# exactly like typing them into Excel with a text quote-prefix.
$ws.Range("D2").Value = "31.543.18"
$ws.Range("D3").Value = "1.993.20"
$ws.Range("E3").Value = "  +6.53%  "
$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'0.8125"
$ws.Range("E5").Value = "  +73.16%  "
$ws.Range("D6").Value = "'254.19"
$ws.Range("E6").Value = "  +4.29%  "
$ws.Range("D7").Value = "'0.9995"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.3430"
$ws.Range("E8").Value = "  +19.25%  "
$ws.Range("E9").Value = "  +17.31%  "
$ws.Range("D10").Value = "'0.07097"
$ws.Range("E10").Value = "  +10.13%  "
$ws.Range("D11").Value = "'0.8477"
$ws.Range("E11").Value = "  +17.68%  "
$ws.Range("D12").Value = "'0.08121"
$ws.Range("E12").Value = "  +4.57%  "
$ws.Range("D13").Value = "'102.53"
$ws.Range("E13").Value = "  +7.14%  "
$ws.Range("D14").Value = "1.991.63"
$ws.Range("E14").Value = "  +6.38%  "
$ws.Range("D15").Value = "'5.499"
$ws.Range("E15").Value = "  +7.28%  "
$ws.Range("D16").Value = "'276.43"
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("D17").Value = "31.543.95"
$ws.Range("E17").Value = "  +4.08%  "
$ws.Range("D18").Value = "'14.03"
$ws.Range("E18").Value = "  +8.15%  "
$ws.Range("D19").Value = "'0.000008012"
$ws.Range("E19").Value = "  +7.91%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'5.731"
$ws.Range("E20").Value = "  +9.73%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.250.06"
$ws.Range("E21").Value = "  +6.18%  "
$ws.Range("D22").Value = "'0.9994"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").Value = "'0.9995"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").Value = "'6.937"
$ws.Range("E24").Value = "  +11.39%  "
$ws.Range("D25").Value = "'9.711"
$ws.Range("E25").Value = "  +7.40%  "
$ws.Range("D26").Value = "'0.1532"
$ws.Range("E26").Value = "  +60.26%  "
$ws.Range("D27").Value = "'165.77"
$ws.Range("E27").Value = "  +1.57%  "
$ws.Range("D28").Value = "'19.85"
$ws.Range("E28").Value = "  +6.43%  "
$ws.Range("D29").Value = "'2.215"
$ws.Range("E29").Value = "  +17.94%  "
$ws.Range("D30").Value = "'1.572"
$ws.Range("E30").Value = "  +7.02%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.594"
$ws.Range("E31").Value = "  +8.93%  "
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "'1.351"
$ws.Range("E32").Value = "  +2.43%  "
$ws.Range("D33").Value = "'4.341"
$ws.Range("E33").Value = "  +5.97%  "
$ws.Range("D34").Value = "'0.05203"
$ws.Range("E34").Value = "  +8.29%  "
$ws.Range("D35").Value = "'1.222"
$ws.Range("E35").Value = "  +9.26%  "
$ws.Range("D36").Value = "'0.7555"
$ws.Range("E36").Value = "  +9.99%  "
$ws.Range("D37").Value = "'2.787"
$ws.Range("E37").Value = "  +2.97%  "
$ws.Range("B38").Value = "Frax"
$ws.Range("C38").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D38").Value = "'0.9987"
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.02014"
$ws.Range("E39").Value = "  +7.72%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.940"
$ws.Range("E40").Value = "  +4.61%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'6.661"
$ws.Range("E41").Value = "  +6.96%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'78.68"
$ws.Range("E42").Value = "  +6.14%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.4697"
$ws.Range("E43").Value = "  +10.99%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'2.086"
$ws.Range("E44").Value = "  +8.08%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'107.08"
$ws.Range("E45").Value = "  +6.34%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "'0.8625"
$ws.Range("E46").Value = "  +4.47%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "'1.000"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.998"
$ws.Range("E48").Value = "  +4.65%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "'7.540"
$ws.Range("E49").Value = "  +9.11%  "
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").Value = "'0.4297"
$ws.Range("E50").Value = "  +9.96%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "'36.60"
$ws.Range("E51").Value = "  +4.33%  "
